$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (existing cells, style already present)
$ws.Range('C4').Value = 'VerificaChiusura'
$ws.Range('C5').Value = 'VerificaChiusura'
$ws.Range('A6').Value = 'ComunicazioneEsito'
$ws.Range('B6').Value = 'AccordoRaggiunto-o-NonRichiesto'
$ws.Range('A7').Value = 'InvioEmailEsito'
$ws.Range('B7').Value = 'EmailEsitoInviata'
$ws.Range('C7').Value = 'AttesaRicezioneAtto'
$ws.Range('A8').Value = 'InvioAtto'
$ws.Range('B8').Value = 'AttoInviato'
$ws.Range('C8').Value = 'AttesaRicezioneAtto'
$ws.Range('A9').Value = 'AttesaRicezioneAtto'
$ws.Range('B9').Value = 'ConciliazioneAvviata'
$ws.Range('C9').Value = 'ConciliazionePerizia'
$ws.Range('A10').Value = 'AttesaRicezioneAtto'
$ws.Range('B10').Value = 'AttoFirmato'
$ws.Range('C10').Value = 'VerificaChiusura'
$ws.Range('A11').Value = 'AttesaRicezioneAtto'
$ws.Range('B11').Value = 'ForzaturaChiusura'
$ws.Range('C11').Value = 'VerificaChiusura'
$ws.Range('A12').Value = 'VerificaChiusura'
$ws.Range('B12').Value = 'PeriziaIncompleta'
$ws.Range('C12').Value = 'EsecuzionePerizia'
$ws.Range('A13').Value = 'VerificaChiusura'
$ws.Range('B13').Value = 'DatiObbligatoriMancanti'
$ws.Range('C13').Value = 'VerificaChiusura'
$ws.Range('A14').Value = 'VerificaChiusura'
$ws.Range('B14').Value = 'ChiusuraVerificata'
$ws.Range('C14').Value = 'ChiusuraAutomatica'
$ws.Range('A15').Value = 'VerificaChiusura'
$ws.Range('B15').Value = 'ChiusuraVerificata'
$ws.Range('C15').Value = 'ChiusuraManuale'
$ws.Range('A16').Value = 'GestioneRiapertura'
$ws.Range('B16').Value = 'PeriziaIntegrazione'
$ws.Range('C16').Value = 'EsecuzionePerizia'
$ws.Range('A17').Value = 'GestioneRiapertura'
$ws.Range('C17').Value = 'GestioneContestazione'

# Set new cell values, then apply the thin-border style used by the rest of the data rows
$ws.Range('B4').Value = 'Chiusura Post SelfCare'
$ws.Range('B5').Value = 'ChiusuraAvviata'
$ws.Range('B17').Value = 'GestioneContestazione'
$ws.Range('A18').Value = 'GestioneRiapertura'
$ws.Range('B18').Value = 'NonRiapertura'
$ws.Range('C18').Value = 'SceltaManualeServizio'
$ws.Range('A19').Value = 'GestioneContestazione'
$ws.Range('B19').Value = 'PeriziaContestazione'
$ws.Range('C19').Value = 'PeriziaContestazione'
$ws.Range('A20').Value = 'GestioneContestazione'
$ws.Range('B20').Value = 'ContestazioneDaFinalizzare'
$ws.Range('C20').Value = 'FinalizzaContestazione'
$ws.Range('A21').Value = 'FinalizzaContestazione'
$ws.Range('B21').Value = 'ChiusuraAvviata'
$ws.Range('C21').Value = 'VerificaChiusura'

# Apply borders to newly created cells so they match style index 2 (thin border all around)
$newCells = @('B4','B5','B17','A18','B18','C18','A19','B19','C19','A20','B20','C20','A21','B21','C21')
foreach ($addr in $newCells) {
    $ws.Range($addr).Borders.LineStyle = 1
}
